$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Rename sheets
$ws1.Name = "ETS_1_monthly"
$ws2.Name = "ETS_1_overall"

# 2. Update header labels with "|K" suffix (temperature unit annotation)
$ws1.Range("C2").Value = "TConEnt_avg|K"
$ws1.Range("I2").Value = "TConEnt_avg|K"
$ws1.Range("O2").Value = "TConEnt_avg|K"
$ws1.Range("D2").Value = "TConLvg_avg|K"
$ws1.Range("J2").Value = "TConLvg_avg|K"
$ws1.Range("P2").Value = "TConLvg_avg|K"
$ws1.Range("E2").Value = "TEvaEnt_avg|K"
$ws1.Range("K2").Value = "TEvaEnt_avg|K"
$ws1.Range("Q2").Value = "TEvaEnt_avg|K"
$ws1.Range("F2").Value = "TEvaLvg_avg|K"
$ws1.Range("L2").Value = "TEvaLvg_avg|K"
$ws1.Range("R2").Value = "TEvaLvg_avg|K"

$ws2.Range("C1").Value = "TEvaEnt_avg|K"
$ws2.Range("D1").Value = "TEvaLvg_avg|K"
$ws2.Range("E1").Value = "TConEnt_avg|K"
$ws2.Range("F1").Value = "TConLvg_avg|K"

# 3. Update sheet2 B1 header to COP_mon (model now reports monthly-style COP label)
$ws2.Range("B1").Value = "COP_mon"

# 4. Update sheet1 numeric grid (rows 4-15), clearing cells no longer populated
# Row 4
$ws1.Range("B4").Value = 7.942139933204688
$ws1.Range("C4").Value = 304.3406020454739
$ws1.Range("D4").Value = 306.6870720904806
$ws1.Range("E4").Value = 284.7343922490659
$ws1.Range("F4").Value = 282.6839639414912
$ws1.Range("G4").Value = 184
$ws1.Range("H4").Value = 4.577358531884396
$ws1.Range("I4").Value = 313.0036964416504
$ws1.Range("J4").Value = 321.1116828918457
$ws1.Range("K4").Value = 285.4331016540527
$ws1.Range("L4").Value = 279.0946235656738
$ws1.Range("M4").Value = 16
$ws1.Range("N4").ClearContents()
$ws1.Range("O4").ClearContents()
$ws1.Range("P4").ClearContents()
$ws1.Range("Q4").ClearContents()
$ws1.Range("R4").ClearContents()
$ws1.Range("S4").ClearContents()

# Row 5
$ws1.Range("B5").Value = 7.398032102648034
$ws1.Range("C5").Value = 304.2139945194639
$ws1.Range("D5").Value = 307.9323849382072
$ws1.Range("E5").Value = 284.5940964271282
$ws1.Range("F5").Value = 281.3789355047818
$ws1.Range("G5").Value = 580
$ws1.Range("H5").ClearContents()
$ws1.Range("I5").ClearContents()
$ws1.Range("J5").ClearContents()
$ws1.Range("K5").ClearContents()
$ws1.Range("L5").ClearContents()
$ws1.Range("M5").ClearContents()
$ws1.Range("N5").ClearContents()
$ws1.Range("O5").ClearContents()
$ws1.Range("P5").ClearContents()
$ws1.Range("Q5").ClearContents()
$ws1.Range("R5").ClearContents()
$ws1.Range("S5").ClearContents()

# Row 6
$ws1.Range("B6").Value = 7.178666154833804
$ws1.Range("C6").Value = 307.9058723449707
$ws1.Range("D6").Value = 310.3665981292725
$ws1.Range("E6").Value = 287.2326602935791
$ws1.Range("F6").Value = 285.1433849334717
$ws1.Range("G6").Value = 32
$ws1.Range("H6").Value = 3.74583039211916
$ws1.Range("I6").Value = 320.0131688271799
$ws1.Range("J6").Value = 327.6829728977655
$ws1.Range("K6").Value = 285.7267282137307
$ws1.Range("L6").Value = 280.1015549526419
$ws1.Range("M6").Value = 93
$ws1.Range("N6").ClearContents()
$ws1.Range("O6").ClearContents()
$ws1.Range("P6").ClearContents()
$ws1.Range("Q6").ClearContents()
$ws1.Range("R6").ClearContents()
$ws1.Range("S6").ClearContents()

# Row 7
$ws1.Range("B7").Value = 6.603569373483853
$ws1.Range("C7").Value = 309.4853820800781
$ws1.Range("D7").Value = 311.9534301757812
$ws1.Range("E7").Value = 287.2582528250558
$ws1.Range("F7").Value = 285.1903337751116
$ws1.Range("G7").Value = 14
$ws1.Range("H7").Value = 3.732659283004948
$ws1.Range("I7").Value = 319.9744325295473
$ws1.Range("J7").Value = 327.3189947666266
$ws1.Range("K7").Value = 285.6676283616286
$ws1.Range("L7").Value = 280.2876970340044
$ws1.Range("M7").Value = 78
$ws1.Range("N7").ClearContents()
$ws1.Range("O7").ClearContents()
$ws1.Range("P7").ClearContents()
$ws1.Range("Q7").ClearContents()
$ws1.Range("R7").ClearContents()
$ws1.Range("S7").ClearContents()

# Row 8
$ws1.Range("B8").Value = 7.28299541021487
$ws1.Range("C8").Value = 307.3153228759766
$ws1.Range("D8").Value = 309.5024375915527
$ws1.Range("E8").Value = 286.9843673706055
$ws1.Range("F8").Value = 285.1379203796387
$ws1.Range("G8").Value = 16
$ws1.Range("H8").Value = 3.577547628486173
$ws1.Range("I8").Value = 321.2474238557635
$ws1.Range("J8").Value = 329.1511213194649
$ws1.Range("K8").Value = 285.5648481261055
$ws1.Range("L8").Value = 279.8674385502653
$ws1.Range("M8").Value = 106
$ws1.Range("N8").ClearContents()
$ws1.Range("O8").ClearContents()
$ws1.Range("P8").ClearContents()
$ws1.Range("Q8").ClearContents()
$ws1.Range("R8").ClearContents()
$ws1.Range("S8").ClearContents()

# Row 9
$ws1.Range("B9").Value = 7.298668928368913
$ws1.Range("C9").Value = 304.2545438228169
$ws1.Range("D9").Value = 308.2897540899935
$ws1.Range("E9").Value = 284.4730277775888
$ws1.Range("F9").Value = 280.9909352159666
$ws1.Range("G9").Value = 574
$ws1.Range("H9").ClearContents()
$ws1.Range("I9").ClearContents()
$ws1.Range("J9").ClearContents()
$ws1.Range("K9").ClearContents()
$ws1.Range("L9").ClearContents()
$ws1.Range("M9").ClearContents()
$ws1.Range("N9").ClearContents()
$ws1.Range("O9").ClearContents()
$ws1.Range("P9").ClearContents()
$ws1.Range("Q9").ClearContents()
$ws1.Range("R9").ClearContents()
$ws1.Range("S9").ClearContents()

# Row 10
$ws1.Range("B10").Value = 7.495356099651389
$ws1.Range("C10").Value = 304.2288258111299
$ws1.Range("D10").Value = 307.6599715503294
$ws1.Range("E10").Value = 284.5827935062238
$ws1.Range("F10").Value = 281.6092470083664
$ws1.Range("G10").Value = 536
$ws1.Range("H10").ClearContents()
$ws1.Range("I10").ClearContents()
$ws1.Range("J10").ClearContents()
$ws1.Range("K10").ClearContents()
$ws1.Range("L10").ClearContents()
$ws1.Range("M10").ClearContents()
$ws1.Range("N10").ClearContents()
$ws1.Range("O10").ClearContents()
$ws1.Range("P10").ClearContents()
$ws1.Range("Q10").ClearContents()
$ws1.Range("R10").ClearContents()
$ws1.Range("S10").ClearContents()

# Row 11
$ws1.Range("B11").Value = 7.398264981324923
$ws1.Range("C11").Value = 307.3450317382812
$ws1.Range("D11").Value = 309.6193695068359
$ws1.Range("E11").Value = 287.0662892659505
$ws1.Range("F11").Value = 285.1390635172526
$ws1.Range("G11").Value = 12
$ws1.Range("H11").Value = 3.903475843913796
$ws1.Range("I11").Value = 318.3003565470378
$ws1.Range("J11").Value = 325.9420687357585
$ws1.Range("K11").Value = 285.6641527811686
$ws1.Range("L11").Value = 279.9779968261719
$ws1.Range("M11").Value = 48
$ws1.Range("N11").Value = 5.090269250889454
$ws1.Range("O11").Value = 310.7619323730469
$ws1.Range("P11").Value = 319.2355041503906
$ws1.Range("Q11").Value = 287.1790466308594
$ws1.Range("R11").Value = 280.3726196289062
$ws1.Range("S11").Value = 4

# Row 12
$ws1.Range("B12").Value = 7.801118574982373
$ws1.Range("C12").Value = 304.2833319223055
$ws1.Range("D12").Value = 307.0497262606057
$ws1.Range("E12").Value = 284.9828753727738
$ws1.Range("F12").Value = 282.5713207285891
$ws1.Range("G12").Value = 372
$ws1.Range("H12").Value = 4.941688848859474
$ws1.Range("I12").Value = 310.209716796875
$ws1.Range("J12").Value = 319.3551025390625
$ws1.Range("K12").Value = 285.4824066162109
$ws1.Range("L12").Value = 278.1847839355469
$ws1.Range("M12").Value = 4
$ws1.Range("N12").ClearContents()
$ws1.Range("O12").ClearContents()
$ws1.Range("P12").ClearContents()
$ws1.Range("Q12").ClearContents()
$ws1.Range("R12").ClearContents()
$ws1.Range("S12").ClearContents()

# Row 13
$ws1.Range("B13").Value = 8.092016102171307
$ws1.Range("C13").Value = 304.6202873461174
$ws1.Range("D13").Value = 306.7063894560843
$ws1.Range("E13").Value = 285.4050653631037
$ws1.Range("F13").Value = 283.580605246804
$ws1.Range("G13").Value = 66
$ws1.Range("H13").Value = 4.106205733069801
$ws1.Range("I13").Value = 316.8456132676866
$ws1.Range("J13").Value = 325.5143636067708
$ws1.Range("K13").Value = 286.3556230333116
$ws1.Range("L13").Value = 279.7928449842665
$ws1.Range("M13").Value = 36
$ws1.Range("N13").Value = 4.390349511591186
$ws1.Range("O13").Value = 315.6627960205078
$ws1.Range("P13").Value = 325.0311431884766
$ws1.Range("Q13").Value = 287.5096435546875
$ws1.Range("R13").Value = 280.2787322998047
$ws1.Range("S13").Value = 4

# Row 14
$ws1.Range("B14").Value = 8.022839834009444
$ws1.Range("C14").Value = 304.4334746979095
$ws1.Range("D14").Value = 306.843497140067
$ws1.Range("E14").Value = 285.469965002039
$ws1.Range("F14").Value = 283.362056857937
$ws1.Range("G14").Value = 182
$ws1.Range("H14").Value = 4.621800966566261
$ws1.Range("I14").Value = 313.1015973772322
$ws1.Range("J14").Value = 322.3745727539062
$ws1.Range("K14").Value = 286.5446428571428
$ws1.Range("L14").Value = 279.2734069824219
$ws1.Range("M14").Value = 14
$ws1.Range("N14").ClearContents()
$ws1.Range("O14").ClearContents()
$ws1.Range("P14").ClearContents()
$ws1.Range("Q14").ClearContents()
$ws1.Range("R14").ClearContents()
$ws1.Range("S14").ClearContents()

# Row 15
$ws1.Range("B15").Value = 7.631717645535224
$ws1.Range("C15").Value = 304.2040211465122
$ws1.Range("D15").Value = 307.3775184956772
$ws1.Range("E15").Value = 284.7950769216528
$ws1.Range("F15").Value = 282.0382810475137
$ws1.Range("G15").Value = 422
$ws1.Range("H15").Value = 7.728417475906985
$ws1.Range("I15").Value = 307.4984741210937
$ws1.Range("J15").Value = 310.1500244140625
$ws1.Range("K15").Value = 289.1929016113281
$ws1.Range("L15").Value = 286.8865051269531
$ws1.Range("M15").Value = 2
$ws1.Range("N15").ClearContents()
$ws1.Range("O15").ClearContents()
$ws1.Range("P15").ClearContents()
$ws1.Range("Q15").ClearContents()
$ws1.Range("R15").ClearContents()
$ws1.Range("S15").ClearContents()

# 5. Update sheet2 numeric grid (rows 2-4)
# Row 2
$ws2.Range("B2").Value = 4.697357491249054
$ws2.Range("C2").Value = 287.3443450927734
$ws2.Range("D2").Value = 280.3256759643555
$ws2.Range("E2").Value = 313.2123641967773
$ws2.Range("F2").Value = 322.1333236694336
$ws2.Range("G2").Value = 8

# Row 3
$ws2.Range("B3").Value = 7.516163328096475
$ws2.Range("C3").Value = 284.7888379355338
$ws2.Range("D3").Value = 281.9301049682209
$ws2.Range("E3").Value = 304.355131423593
$ws2.Range("F3").Value = 307.6541029710036
$ws2.Range("G3").Value = 2990

# Row 4
$ws2.Range("B4").Value = 3.820315707019815
$ws2.Range("C4").Value = 285.753366388662
$ws2.Range("D4").Value = 279.9777671371959
$ws2.Range("E4").Value = 319.1527301778721
$ws2.Range("F4").Value = 326.9720522018164
$ws2.Range("G4").Value = 397

